# Daily attendance processing - 2025-10-20 09:23:07
# Updates "Recorded By" (system/backup account labels), a few Students
# fraction counts, and a few Average Attendance % figures on the
# attendance-analysis worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Helper: assign a plain-text value to a cell. Percentage-looking text
# (e.g. "72.3%") would otherwise be auto-converted by Excel into a
# numeric percentage (changing both the stored type and the cell's
# style index), so for those we briefly force literal-text entry with a
# leading apostrophe, then restore the cell's original number
# format/style by pasting the formats (only) from a same-styled
# neighbour cell that this script leaves untouched.
# ---------------------------------------------------------------------
function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text,
        [string]$FormatSourceAddress
    )

    if ($FormatSourceAddress) {
        $ws.Range($Address).Value = "'" + $Text
        $ws.Range($FormatSourceAddress).Copy() | Out-Null
        $ws.Range($Address).PasteSpecial(-4122) | Out-Null
    } else {
        $ws.Range($Address).Value = $Text
    }
}

# ----------------------- "Recorded By" (column G) -----------------------
Set-TextValue "G2"  "system, backup@backdoor.com, System"
Set-TextValue "G4"  "backup@backdoor.com, System"
Set-TextValue "G5"  "backup@backdoor.com, System"
Set-TextValue "G8"  "backup@backdoor.com, System"
Set-TextValue "G22" "dnasr281@gmail.com, System"
Set-TextValue "G29" "system, backup@backdoor.com, System"
Set-TextValue "G31" "backup@backdoor.com, System"
Set-TextValue "G32" "backup@backdoor.com, System"
Set-TextValue "G35" "backup@backdoor.com, System"
Set-TextValue "G49" "dnasr281@gmail.com, System"
Set-TextValue "G56" "system, backup@backdoor.com, System"
Set-TextValue "G58" "backup@backdoor.com, System"
Set-TextValue "G59" "backup@backdoor.com, System"
Set-TextValue "G62" "backup@backdoor.com, System"
Set-TextValue "G76" "dnasr281@gmail.com, System"
Set-TextValue "G83" "backup@backdoor.com, System"
Set-TextValue "G84" "backup@backdoor.com, System"
Set-TextValue "G85" "backup@backdoor.com, System"
Set-TextValue "G109" "backup@backdoor.com, System"
Set-TextValue "G110" "backup@backdoor.com, System"
Set-TextValue "G111" "backup@backdoor.com, System"
Set-TextValue "G135" "backup@backdoor.com, System"
Set-TextValue "G136" "backup@backdoor.com, System"
Set-TextValue "G137" "backup@backdoor.com, System"

# ----------------------- Students counts (column H) -----------------------
Set-TextValue "H48" "41/57"
Set-TextValue "H49" "33/57"
Set-TextValue "H62" "34/55"
Set-TextValue "H75" "42/55"
Set-TextValue "H76" "36/55"

# ------------------- Average Attendance % (columns L, S) -------------------
# These cells carry percentage-formatted text; use an untouched sibling
# cell with the same style (s="4") as the paste-format source so the
# original style index is preserved exactly.
Set-TextValue "L10" "72.3%" "L18"
Set-TextValue "S16" "67.0%" "S18"
Set-TextValue "S17" "62.2%" "S18"
